$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5: add Details text (row 5, Santana et al paper)
$ws.Range("E5").Value = 'Our study showed mixed effects of long-term conservation investment in N2000 farmland. We found positive effects on flagship species, and on species associated with fallows, which were the main targets of conservation investment. Finally, long-term evaluations of conservation investment are required, in order to monitor and improve the effectiveness of billions of euros needed annually for managing N2000. '

# ================= Row 12 =================
$ws.Rows.Item(12).RowHeight = 72.6
$ws.Range("A12").Value = 'Fryxell et al'
$ws.Range("B12").Value = 'Resource management cycles and the sustainability of harvested wildlife populations'
$ws.Range("C12").Value = 2010
$ws.Range("D12").Value = 'Modelled complex dynamic relatinoships between harvest quotas, users, managers, and population survival'
$ws.Range("E12").Value = 'Here we show that weak compensatory response by harvesters or resource managers can itself generate cyclic variation in resources, exacerbating the risk of collapse. Weak harvest regulation contributes to the problem rather than providing an acceptable management solution to resource fluctuation. Our simulations suggest that the risk of population collapse could be dramatically higher in systems with dynamic effort and quota levels (Fig. 3), simply because of extreme population excursions caused by quasiperiodic dynamics resulting from even mild levels of environmental stochasticity.'

$ws.Range("B4").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 13 =================
$ws.Rows.Item(13).RowHeight = 60.6
$ws.Range("A13").Value = 'Armsworth et al'
$ws.Range("B13").Value = 'Is conservation right to go big? PA size and conservation return on invesment'
$ws.Range("C13").Value = 2018
$ws.Range("D13").Value = 'Examine how PA size influences conservation return on investment'
$ws.Range("E13").Value = 'Policy guidelines for creating new protected areas commonly recommend larger protected areas be favored. We examine whether these recommendations are justified, providing the first evaluation of this question to use return-on-investment (ROI) methods that account for how protected area size influences multiple ecological benefits and the economic costs of protection. A portfolio of site sizes may need to be included in protected area networks when multiple objectives motivate conservation.'

$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 14 =================
$ws.Rows.Item(14).RowHeight = 36.6
$ws.Range("A14").Value = 'Lindsey et al '
$ws.Range("B14").Value = 'Underperformance of African PA networks and the case for new conservation models: Insights from Zambia'
$ws.Range("C14").Value = 2014
$ws.Range("D14").Value = 'Assess why PAs in Zambia are not performing well'
$ws.Range("E14").Value = 'They mention increasing human population increases in and around PAs is having a negative effect'

$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 15 =================
$ws.Rows.Item(15).RowHeight = 36.6
$ws.Range("A15").Value = 'Wittemyer et al'
$ws.Range("B15").Value = 'Accelerated human population growth at protected area edges'
$ws.Range("C15").Value = 2008
$ws.Range("D15").Value = 'Assess human pop growth around PAs in 45 countries'
$ws.Range("E15").Value = 'Contrary to predictions of this argument, we found that average human population growth rates on the borders of 306 PAs in 45 countries in Africa and Latin America were nearly double average rural growth, suggesting that PAs attract, rather than repel, human settlement. Human growth rate around PAs correlated with forest loss. '

$ws.Range("B4").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 16 =================
$ws.Rows.Item(16).RowHeight = 24.6
$ws.Range("A16").Value = 'Bruner et al '
$ws.Range("B16").Value = 'Financial costs and shortfalls for expanding PA systems in developing countries'
$ws.Range("C16").Value = 2004
$ws.Range("E16").Value = 'Reference for the efficent use of PA resources and investment '

$ws.Range("B4").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 17 =================
$ws.Rows.Item(17).RowHeight = 72.6
$ws.Range("A17").Value = 'Utami et al'
$ws.Range("B17").Value = 'Prioritizing management strategies to achieve multiple outcomes in a globally significant Indonesian protected area'
$ws.Range("C17").Value = 2020
$ws.Range("D17").Value = 'Assessed different mgmt strategies for cost effectiveness and ability to improve a number of PA values over 15 years'
$ws.Range("E17").Value = 'in this study we aimed to: (a) build an approach capable of assessing the cost, relative bene- fits and cost-effectiveness of implementing threat management strategies that improve a broad range of values in multifunctional areas; (b) bring together and build key information to help managers and stakeholders understand the values, goals, threats, total management costs and opportunities for achieving goals for values, using the TNBB as a case study; and (c) deliver a set of costed, prioritized strategies for achieving goals across multiple important values of the TNBB.'

$ws.Range("A4").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 18 =================
$ws.Rows.Item(18).RowHeight = 24.6
$ws.Range("A18").Value = 'Cullen'
$ws.Range("B18").Value = 'Biodiversity protection prioiritisation: a 25 year review'
$ws.Range("C18").Value = 2012

$ws.Range("A4").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ================= Row 19 =================
$ws.Rows.Item(19).RowHeight = 48.6
$ws.Range("A19").Value = 'Meir et al'
$ws.Range("B19").Value = 'Does conservation planning matter in a dynamic and uncertain world?'
$ws.Range("C19").Value = 2004
$ws.Range("E19").Value = 'Here we explicitly consider the implications for biodiversity conservation of several key assumptions underlying systematic conservation planning methods. Our results suggest that relatively simple rules for deciding which areas to protect outperform both ad hoc investment strategies and comprehensive conservation plans (Figs 1 and 2). This is especially true when degradation rates and uncertainty are high '
$ws.Range("D19").Value = 'Simulate site selection under conditions of different budgets, site acquisition uncertainty'

$ws.Range("A4").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Final selection to match target sheet view
$ws.Activate()
$ws.Range("E25").Select()